# Weekly price-sheet update: insert the new week's record as a new row
# right above the current row 661, pushing the existing rows (old 661..699)
# down by one position (to 662..700).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row at position 661 - everything below shifts down by one.
$ws.Rows(661).Insert()

# Populate the newly inserted row with this week's data.
$ws.Cells.Item(661, 1).Value  = 9
$ws.Cells.Item(661, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(661, 3).Value  = "Metropolitana"
$ws.Cells.Item(661, 4).Value  = 45267
$ws.Cells.Item(661, 5).Value  = 13
$ws.Cells.Item(661, 6).Value  = 100112012
$ws.Cells.Item(661, 7).Value  = "Espinaca"
$ws.Cells.Item(661, 8).Value  = "Sin especificar"
$ws.Cells.Item(661, 9).Value  = "Primera"
$ws.Cells.Item(661, 10).Value = 160
$ws.Cells.Item(661, 11).Value = 8000
$ws.Cells.Item(661, 12).Value = 10000
$ws.Cells.Item(661, 13).Value = 9000
$ws.Cells.Item(661, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(661, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(661, 16).Value = 900
$ws.Cells.Item(661, 17).Value = 10
$ws.Cells.Item(661, 18).Value = "Hortaliza"
